# menambahkan export xlsx pada stock Out
# Adds "No" and "Tanggal Keluar" columns at the front of the DataStok sheet,
# renames several headers, fills in the previously-blank Kode value,
# and appends a new stock-out row (row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift existing data two columns to the right, making room for the
#     new "No" and "Tanggal Keluar" columns -------------------------------
$ws.Range("A1:B1").EntireColumn.Insert()

# --- Header row -----------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "No"
$ws.Cells.Item(1, 2).Value = "Tanggal Keluar"
$ws.Cells.Item(1, 3).Value = "Lokasi"
$ws.Cells.Item(1, 4).Value = "Kode"
$ws.Cells.Item(1, 5).Value = "Nama Barang"
$ws.Cells.Item(1, 6).Value = "Unit"
$ws.Cells.Item(1, 7).Value = "Barang Keluar"
$ws.Cells.Item(1, 8).Value = "Total Barang"
$ws.Cells.Item(1, 9).Value = "ID"

# The "Tanggal Keluar" column (B) holds date/time serials; give it a
# date-time display format (maps to builtin numFmtId 22: m/d/yy h:mm).
$ws.Range("B2:B8").NumberFormat = "m/d/yy h:mm"

# --- Existing data rows: add sequence number + timestamp, fill Kode -------
# "Kode" (column D) is the numeric-looking code "12345" stored as text, so
# it is entered with a leading apostrophe, same as a user typing it in.
$tanggalKeluar = @(
    45242.950240115744,
    45242.95028357639,
    45242.950316724535,
    45242.95065928241,
    45242.950700381945,
    45243.25965025463
)

for ($i = 0; $i -lt $tanggalKeluar.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $tanggalKeluar[$i]
    $ws.Cells.Item($row, 4).Value = "'12345"
}

# --- New stock-out row (row 8) --------------------------------------------
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 45243.848213032405
$ws.Cells.Item(8, 3).Value = "Tangerang"
$ws.Cells.Item(8, 4).Value = "'12345"
$ws.Cells.Item(8, 5).Value = "sukasuka"
$ws.Cells.Item(8, 6).Value = "dus"
$ws.Cells.Item(8, 7).Value = 10000
$ws.Cells.Item(8, 8).Value = 310000
$ws.Cells.Item(8, 9).Value = 2
